# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#    populated with the quarterly fund-holdings table (mirrors the layout
#    of "2021-Q4").
# 2. Insert a new leading data row into "总计" for the "2022-Q1" totals,
#    pushing the existing "2021-Q4" row down and renumbering its index.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item(1)      # "2021-Q4" - used as a formatting template

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: worksheet references resolve by current tab position, so any
# handle obtained *before* this insert (e.g. to "总计") would silently
# point at the wrong tab afterwards. Fetch it fresh, now that the sheet
# order is final.
$wsTotal = $wb.Worksheets.Item("总计")

# -- header row -----------------------------------------------------
$headers = @{
    "B1" = "基金代码"
    "C1" = "基金名称"
    "D1" = "基金规模"
    "E1" = "股票总仓位"
    "F1" = "仓位占比"
    "G1" = "持有市值(亿元)"
    "H1" = "仓位排名"
}
foreach ($addr in $headers.Keys) {
    $wsQ1.Range($addr).Value = $headers[$addr]
}

# -- data rows --------------------------------------------------------
# columns B-G are free-form text in the source data (fund codes such as
# "007835" must keep their leading zero, and figures such as "3.00" must
# keep their trailing zero) so force the Text number format before typing
# the values in, just like typing into a pre-formatted "Text" column.
$textRange = $wsQ1.Range("B2:G4")
$textRange.NumberFormat = "@"

$rows = @(
    @{ A=0; B="320003"; C="诺安先锋混合";             D="45.79"; E="69.96"; F="3.00"; G="1.3737"; H=8  },
    @{ A=1; B="007835"; C="国泰鑫睿混合";             D="9.37";  E="78.94"; F="2.73"; G="0.2558"; H=10 },
    @{ A=2; B="001743"; C="诺安优选回报灵活配置混合"; D="6.13";  E="71.32"; F="3.73"; G="0.2286"; H=9  }
)

$r = 2
foreach ($row in $rows) {
    $wsQ1.Range("A$r").Value = $row.A
    $wsQ1.Range("B$r").Value = $row.B
    $wsQ1.Range("C$r").Value = $row.C
    $wsQ1.Range("D$r").Value = $row.D
    $wsQ1.Range("E$r").Value = $row.E
    $wsQ1.Range("F$r").Value = $row.F
    $wsQ1.Range("G$r").Value = $row.G
    $wsQ1.Range("H$r").Value = $row.H
    $r = $r + 1
}

# -- formatting: reuse the bold/bordered header & index style already
#    used on the other sheets (style of "2021-Q4"!B1 / "2021-Q4"!A2) ---
$wsQ4.Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. "总计" sheet: insert a new leading row for the "2022-Q1" totals
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 1.86

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1

# restore the index-column style on the new row (insert leaves it
# without the bold/bordered formatting applied elsewhere in column A)
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# the row that shifted down keeps the original body-cell formatting, but
# the inserted blank cells may have inherited stray formatting - clear
# it so B2:D2 match the plain (unstyled) body cells used elsewhere
$wsTotal.Range("B2:D2").ClearFormats()

$excel.CutCopyMode = $false

# keep the original active tab ("2021-Q4"), since adding/activating the
# new sheet along the way left it selected
$wsQ4.Activate()
